$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("周报")

# Clear B2 (old text "1.学会eclipse的maven+spring项目的构建。" removed)
$ws.Range("B2").ClearContents()

# C2 now holds what used to be the unique string "能创建数据库表"
$ws.Range("C2").Value = "能创建数据库表"

# D2 now holds the shortened "继续学习" (was "继续学习创建数据库表")
$ws.Range("D2").Value = "继续学习"

# B3 also references the same "能创建数据库表" string
$ws.Range("B3").Value = "能创建数据库表"
